$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the ratio_threshold_range Min value (row 4)
$ws.Range("B4").Value = 0.7

# Delete the theta_threshold_range row (row 5) entirely; shifts pie_threshold_range up to row 5
$ws.Rows("5").Delete()

# Configure the page setup (paper size + orientation) for printing
$ps = $ws.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1

# Set the active selection to match the target state
$ws.Range("B4:C4").Select()
